# Apply the updated NQ/HQ market-price + profit figures captured by the
# scheduled Marilith_Profits runner, sheet by sheet / row by row.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 39
$ws.Range("I5").Value = 39
$ws.Range("K5").Value = 39
$ws.Range("M5").Value = 76
$ws.Range("H12").Value = 553.8
$ws.Range("I12").Value = 621.8570999999999
$ws.Range("J12").Value = 395
$ws.Range("K12").Value = 621.8570999999999
$ws.Range("L12").Value = 395
$ws.Range("M12").Value = -451.8570999999999
$ws.Range("N12").Value = -735
$ws.Range("H38").Value = 291.58334
$ws.Range("I38").Value = 77.666664
$ws.Range("J38").Value = 933.3333
$ws.Range("K38").Value = 232.999992
$ws.Range("L38").Value = 2799.9999
$ws.Range("M38").Value = 139.000008
$ws.Range("N38").Value = -3543.9999
$ws.Range("H40").Value = 4381.364
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 4519.5
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 4519.5
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -4869.5
$ws.Range("H58").Value = 5100
$ws.Range("J58").Value = 5100
$ws.Range("L58").Value = 15300
$ws.Range("N58").Value = -15600
$ws.Range("H98").Value = 2339.9285
$ws.Range("I98").Value = 1206.8
$ws.Range("J98").Value = 2969.4443
$ws.Range("K98").Value = 1206.8
$ws.Range("L98").Value = 2969.4443
$ws.Range("M98").Value = 291.2
$ws.Range("N98").Value = -5965.4443
$ws.Range("H116").Value = 6273.125
$ws.Range("I116").Value = 5530.8335
$ws.Range("K116").Value = 5530.8335
$ws.Range("M116").Value = -2088.8335
$ws.Range("H122").Value = 2339.9285
$ws.Range("I122").Value = 1206.8
$ws.Range("J122").Value = 2969.4443
$ws.Range("K122").Value = 3620.4
$ws.Range("L122").Value = 8908.332900000001
$ws.Range("M122").Value = -1170.4
$ws.Range("N122").Value = -13808.3329
$ws.Range("H138").Value = 3744.2104
$ws.Range("J138").Value = 3744.2104
$ws.Range("L138").Value = 11232.6312
$ws.Range("N138").Value = -21512.6312

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1867.0714
$ws.Range("I74").Value = 1241.6923
$ws.Range("K74").Value = 1241.6923
$ws.Range("M74").Value = -367.6922999999999
$ws.Range("H77").Value = 1867.0714
$ws.Range("I77").Value = 1241.6923
$ws.Range("K77").Value = 6208.461499999999
$ws.Range("M77").Value = -1840.461499999999
$ws.Range("H97").Value = 1224.0769
$ws.Range("I97").Value = 702.4
$ws.Range("K97").Value = 702.4
$ws.Range("M97").Value = -206.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9499.5
$ws.Range("I20").Value = 9499.5
$ws.Range("K20").Value = 9499.5
$ws.Range("M20").Value = -9252.5
$ws.Range("H86").Value = 1903.5
$ws.Range("I86").Value = 1800
$ws.Range("K86").Value = 1800
$ws.Range("M86").Value = -677
$ws.Range("H89").Value = 1903.5
$ws.Range("I89").Value = 1800
$ws.Range("K89").Value = 9000
$ws.Range("M89").Value = -3384
$ws.Range("H94").Value = 1681.45
$ws.Range("I94").Value = 1044.9333
$ws.Range("K94").Value = 1044.9333
$ws.Range("M94").Value = -593.9332999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3337
$ws.Range("I31").Value = 2734.8333
$ws.Range("J31").Value = 3665.4546
$ws.Range("K31").Value = 2734.8333
$ws.Range("L31").Value = 3665.4546
$ws.Range("M31").Value = -2439.8333
$ws.Range("N31").Value = -4255.4546
$ws.Range("H34").Value = 3337
$ws.Range("I34").Value = 2734.8333
$ws.Range("J34").Value = 3665.4546
$ws.Range("K34").Value = 2734.8333
$ws.Range("L34").Value = 3665.4546
$ws.Range("M34").Value = -2532.8333
$ws.Range("N34").Value = -4069.4546
$ws.Range("H58").Value = 2432.25
$ws.Range("I58").Value = 2637.1428
$ws.Range("J58").Value = 998
$ws.Range("K58").Value = 2637.1428
$ws.Range("L58").Value = 998
$ws.Range("M58").Value = -2434.1428
$ws.Range("N58").Value = -1404
$ws.Range("H134").Value = 3906.4666
$ws.Range("I134").Value = 3969.3076
$ws.Range("K134").Value = 11907.9228
$ws.Range("M134").Value = -9372.9228
$ws.Range("H136").Value = 2432.25
$ws.Range("I136").Value = 2637.1428
$ws.Range("J136").Value = 998
$ws.Range("K136").Value = 7911.428400000001
$ws.Range("L136").Value = 2994
$ws.Range("M136").Value = -5361.428400000001
$ws.Range("N136").Value = -8094

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 83.09090999999999
$ws.Range("J2").Value = 52
$ws.Range("K2").Value = 498.5454599999999
$ws.Range("L2").Value = 312
$ws.Range("M2").Value = -385.5454599999999
$ws.Range("N2").Value = -538
$ws.Range("H11").Value = 334749.5
$ws.Range("I11").Value = 666833.3
$ws.Range("K11").Value = 2000499.9
$ws.Range("M11").Value = -2000359.9
$ws.Range("H38").Value = 30.285715
$ws.Range("I38").Value = 27.09091
$ws.Range("J38").Value = 42
$ws.Range("K38").Value = 81.27273
$ws.Range("L38").Value = 126
$ws.Range("M38").Value = 265.72727
$ws.Range("N38").Value = -820
$ws.Range("H107").Value = 2287.25
$ws.Range("I107").Value = 2199
$ws.Range("K107").Value = 6597
$ws.Range("M107").Value = -4677
$ws.Range("H110").Value = 3000
$ws.Range("I110").Value = 3000
$ws.Range("K110").Value = 9000
$ws.Range("M110").Value = -4910
$ws.Range("H121").Value = 2723
$ws.Range("J121").Value = 2903.75
$ws.Range("L121").Value = 8711.25
$ws.Range("N121").Value = -11331.25
$ws.Range("H132").Value = 1348
$ws.Range("I132").Value = 1348
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12132
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9602
$ws.Range("N132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6974.75
$ws.Range("I70").Value = 6450
$ws.Range("J70").Value = 7499.5
$ws.Range("K70").Value = 6450
$ws.Range("L70").Value = 7499.5
$ws.Range("M70").Value = -6180
$ws.Range("N70").Value = -8039.5
$ws.Range("H73").Value = 6974.75
$ws.Range("I73").Value = 6450
$ws.Range("J73").Value = 7499.5
$ws.Range("K73").Value = 6450
$ws.Range("L73").Value = 7499.5
$ws.Range("M73").Value = -5514
$ws.Range("N73").Value = -9371.5
$ws.Range("H113").Value = 762
$ws.Range("I113").Value = 762
$ws.Range("K113").Value = 762
$ws.Range("M113").Value = 1408

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2850
$ws.Range("I7").Value = 2850
$ws.Range("K7").Value = 2850
$ws.Range("M7").Value = -2738
$ws.Range("H40").Value = 5617.273
$ws.Range("I40").Value = 6061.5
$ws.Range("J40").Value = 4432.6665
$ws.Range("K40").Value = 6061.5
$ws.Range("L40").Value = 4432.6665
$ws.Range("M40").Value = -5925.5
$ws.Range("N40").Value = -4704.6665
$ws.Range("H61").Value = 2067.3
$ws.Range("I61").Value = 1891.4445
$ws.Range("J61").Value = 3650
$ws.Range("K61").Value = 1891.4445
$ws.Range("L61").Value = 3650
$ws.Range("M61").Value = -1689.4445
$ws.Range("N61").Value = -4054
$ws.Range("H113").Value = 2067.3
$ws.Range("I113").Value = 1891.4445
$ws.Range("J113").Value = 3650
$ws.Range("K113").Value = 1891.4445
$ws.Range("L113").Value = 3650
$ws.Range("M113").Value = 278.5554999999999
$ws.Range("N113").Value = -7990
$ws.Range("H126").Value = 2850
$ws.Range("I126").Value = 2850
$ws.Range("K126").Value = 8550
$ws.Range("M126").Value = -6080
$ws.Range("H132").Value = 13556.235
$ws.Range("I132").Value = 13490
$ws.Range("K132").Value = 40470
$ws.Range("M132").Value = -37940

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1213.6
$ws.Range("I96").Value = 889.6667
$ws.Range("J96").Value = 1699.5
$ws.Range("K96").Value = 889.6667
$ws.Range("L96").Value = 1699.5
$ws.Range("M96").Value = 483.3333
$ws.Range("N96").Value = -4445.5
$ws.Range("H107").Value = 450
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -5040
$ws.Range("H136").Value = 2436.7646
$ws.Range("I136").Value = 2339.0625
$ws.Range("K136").Value = 7017.1875
$ws.Range("M136").Value = -4467.1875

Write-Output "Applied 217 Marilith_Profits cell updates"
